$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture existing values (columns B:F, rows 1:3) before overwriting anything.
# Use Value2 for reads -- this COM shim's Value *getter* is unreliable, but the
# setter works fine for both Value and Value2.
$header = @($ws.Range("B1").Value2, $ws.Range("C1").Value2, $ws.Range("D1").Value2, $ws.Range("E1").Value2, $ws.Range("F1").Value2)
$row2   = @($ws.Range("B2").Value2, $ws.Range("C2").Value2, $ws.Range("D2").Value2, $ws.Range("E2").Value2, $ws.Range("F2").Value2)
$row3   = @($ws.Range("B3").Value2, $ws.Range("C3").Value2, $ws.Range("D3").Value2, $ws.Range("E3").Value2, $ws.Range("F3").Value2)

# Give the new A1:E1 header range the same format (bold font + border) the old
# header row (B1:F1) already carried, by copying B1's formatting over to A1.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null   # -4122 = xlPasteFormats
$excel.CutCopyMode = $false

# Shift everything one column to the left (B->A, C->B, D->C, E->D, F->E).
# The old column A values (row labels 4 / 14, with their border/bold style) are dropped entirely.
$ws.Range("A1").Value = $header[0]
$ws.Range("B1").Value = $header[1]
$ws.Range("C1").Value = $header[2]
$ws.Range("D1").Value = $header[3]
$ws.Range("E1").Value = $header[4]

$ws.Range("A2").Value = $row2[0]
$ws.Range("B2").Value = $row2[1]
$ws.Range("C2").Value = $row2[2]
$ws.Range("D2").Value = $row2[3]
$ws.Range("E2").Value = $row2[4]

$ws.Range("A3").Value = $row3[0]
$ws.Range("B3").Value = $row3[1]
$ws.Range("C3").Value = $row3[2]
$ws.Range("D3").Value = $row3[3]
$ws.Range("E3").Value = $row3[4]

# Column A on rows 2/3 used to hold the bordered/bold style (old A2/A3); that style must
# not carry over to the new A2/A3 (which now hold plain shifted-in data), so strip it.
$ws.Range("A2:A3").ClearFormats()

# Drop the now-empty former column F.
$ws.Range("F1:F3").Clear()
